# UPDATE removendo a nossa rede da apresentação
#
# Replaces references to "nossa BlockChain" with "rede ethereum" in the
# two slides that talk about the platform's underlying chain technology.
# Uses character-range (TextRange.Characters) edits so only the targeted
# substrings change and the rest of each run/paragraph's formatting
# (err="1" spell-flags, lang, etc.) is left untouched.

$p = $ppt.ActivePresentation

function Replace-InRange {
    param($TextRange, $Find, $Replacement, $Last)
    $current = $TextRange.Text
    if ($Last) {
        $idx = $current.LastIndexOf($Find)
    } else {
        $idx = $current.IndexOf($Find)
    }
    if ($idx -lt 0) {
        return $false
    }
    $sub = $TextRange.Characters($idx + 1, $Find.Length)
    $sub.Text = $Replacement
    return $true
}

# Slide 3 - "Quem Somos?" : "Através da nossa BlockChain" -> "Através da rede ethereum"
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item("Content Placeholder 2")
$tr3 = $shape3.TextFrame.TextRange
Replace-InRange $tr3 "Através da nossa " "Através da rede " $false | Out-Null
Replace-InRange $tr3 "BlockChain" "ethereum" $false | Out-Null

# Slide 7 - "Conceito por trás" : "...submete-lo a nossa BlockChain" -> "...submete-lo a na rede ethereum"
$slide7 = $p.Slides.Item(7)
$shape7 = $slide7.Shapes.Item("Content Placeholder 2")
$tr7 = $shape7.TextFrame.TextRange
Replace-InRange $tr7 " a empresa ou pessoa que deseja realizar o IPO do seu negócio pode clonar nosso contrato e submete-lo a nossa " " a empresa ou pessoa que deseja realizar o IPO do seu negócio pode clonar nosso contrato e submete-lo a na rede " $false | Out-Null
Replace-InRange $tr7 "BlockChain" "ethereum" $true | Out-Null
